$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1148
$ws1.Range("F4").Value = 2601
$ws1.Range("F5").Value = 226

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1148
$ws4.Range("F6").Value = 2601
$ws4.Range("F8").Value = 226
